$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize connector words ("de" -> "De", "el" -> "El", "y" -> "Y", "los" -> "Los")
# in municipality / state names
$ws.Range("B5").Value = "Comitán De Domínguez"
$ws.Range("B7").Value = "Mazapa De Madero"
$ws.Range("A11").Value = "Ciudad De México"
$ws.Range("A22").Value = "Estado De México"
$ws.Range("B22").Value = "Ecatepec De Morelos"
$ws.Range("B24").Value = "Tlalnepantla De Baz"
$ws.Range("B26").Value = "Apaseo El Alto"
$ws.Range("B28").Value = "Acapulco De Juárez"
$ws.Range("B30").Value = "Ayutla De Los Libres"
$ws.Range("B32").Value = "Coyuca De Catalán"
$ws.Range("B33").Value = "Zihuatanejo De Azueta"
$ws.Range("B45").Value = "Pachuca De Soto"
$ws.Range("B46").Value = "Tenango De Doria"
$ws.Range("B48").Value = "Tulancingo De Bravo"
$ws.Range("B50").Value = "Cuautitlán De García Barragán"
$ws.Range("B52").Value = "Unión De Tula"
$ws.Range("B60").Value = "San Nicolás De Los Garza"
$ws.Range("B63").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B72").Value = "Izúcar De Matamoros"
$ws.Range("B76").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B80").Value = "Jalpan De Serra"
$ws.Range("B81").Value = "Landa De Matamoros"
$ws.Range("B104").Value = "Martínez De La Torre"

# Remove footer/notes rows 114-118 (sample size, source, credits, date)
$ws.Range("A114:A118").EntireRow.Delete()
